$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.1909339427948
$ws.Range("B1").Value = 2.371184110641479
$ws.Range("C1").Value = 4.178069114685059
$ws.Range("D1").Value = 2.892116069793701
$ws.Range("E1").Value = 1.123141884803772
